$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.10359999999999
$ws.Range("C21").Value = -13.23120000000001
$ws.Range("C23").Value = -11.9303
$ws.Range("C25").Value = -11.0455
